$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1877
$ws.Range("F4").Value = 876
$ws.Range("F5").Value = 762
$ws.Range("F6").Value = 13302
$ws.Range("F7").Value = 13176
$ws.Range("F9").Value = 771
$ws.Range("F13").Value = 669
$ws.Range("F14").Value = 2090
$ws.Range("F16").Value = 44
$ws.Range("F17").Value = 70
$ws.Range("F19").Value = 394
$ws.Range("F20").Value = 243
$ws.Range("F21").Value = 285
$ws.Range("F22").Value = 415
$ws.Range("F24").Value = 14

# Sheet "演出" (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 27

# Sheet "本地生活" (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 35

# Sheet "全部类型" (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1877
$ws.Range("F5").Value = 876
$ws.Range("F7").Value = 762
$ws.Range("F8").Value = 13302
$ws.Range("F9").Value = 13176
$ws.Range("F11").Value = 771
$ws.Range("F15").Value = 669
$ws.Range("F18").Value = 2090
$ws.Range("F20").Value = 44
$ws.Range("F21").Value = 70
$ws.Range("F25").Value = 35
$ws.Range("F26").Value = 394
$ws.Range("F27").Value = 243
$ws.Range("F28").Value = 285
$ws.Range("F29").Value = 415
$ws.Range("F33").Value = 14
$ws.Range("F34").Value = 27

$wb.Save()
